$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C (Förändrad) for rows 2-15 from 45175 to 45183 (serial date values)
$ws.Range("C2:C15").Value = 45183
